# Applies the addition of columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style used by the
# existing header row (copy format from the adjacent H1 header cell).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(3,8,4,6,1,1,8,4,1,5,1,6,6,5,1,9,10,7,5,7,8,6,7,6,5,8,7,6,9,6,7,7,6,1,6,8,7,6,6,8,4,7,6,5)
$jValues = @(3,8,5,6,1,2,8,5,1,5,1,6,6,5,1,9,10,7,5,7,8,6,7,7,6,8,7,6,9,6,7,7,6,1,6,8,7,6,6,8,4,7,6,5)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
